$d = $word.ActiveDocument

function New-XmlFrag($inner) {
    $pre = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $post = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pre + $inner + $post
}

# 1) "Alberto Perales " -> "Alberto " + proofErr(spellStart) + "Perales" + proofErr(spellEnd) + " "
$p1 = $d.Paragraphs(1).Range
$inner1 = '<w:p><w:r><w:t xml:space="preserve">Alberto </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Perales</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$p1.InsertXML((New-XmlFrag $inner1))

# 2) "b) sub goals are to find right pairing of animals and seed withing pairs in boat "
$p19 = $d.Paragraphs(19).Range
$inner19 = '<w:p><w:r><w:t xml:space="preserve">b) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>sub</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> goals are to find right pairing of animals and seed </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>withing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pairs in boat </w:t></w:r></w:p>'
$p19.InsertXML((New-XmlFrag $inner19))

# 3) "find pairing to travel ….Cat and man , seed and man, parrot and man " (2nd run in paragraph 22)
$p22 = $d.Paragraphs(22).Range
$inner22 = '<w:p><w:r><w:t xml:space="preserve">a) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>find</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> pairing to travel &#8230;.Cat and man , seed and man, parrot and man </w:t></w:r></w:p>'
$p22.InsertXML((New-XmlFrag $inner22))

# 4) "a) found potential solution to pair traveling and they meet goal of not being left with wrong pair "
$p25 = $d.Paragraphs(25).Range
$inner25 = '<w:p><w:r><w:t xml:space="preserve">a) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>found</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> potential solution to pair traveling and they meet goal of not being left with wrong pair </w:t></w:r></w:p>'
$p25.InsertXML((New-XmlFrag $inner25))

# 5) "b) tried to take one at a time but it would leave impossible pairing on either side. "
$p30 = $d.Paragraphs(30).Range
$inner30 = '<w:p><w:r><w:t xml:space="preserve">b) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>tried</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to take one at a time but it would leave impossible pairing on either side. </w:t></w:r></w:p>'
$p30.InsertXML((New-XmlFrag $inner30))

# 6) "2." paragraph (sock problem #2 heading) gains a tab stop + new run
$p39 = $d.Paragraphs(39).Range
$inner39 = '<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2608"/></w:tabs></w:pPr><w:r><w:t>2.</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">The solution for the sub did meet the goal </w:t></w:r></w:p>'
$p39.InsertXML((New-XmlFrag $inner39))

# 7) "a)  " -> "a) " + new run, plus 6 new paragraphs (problem 2 solution content)
$p50 = $d.Paragraphs(50).Range
$inner50 = '<w:p><w:r><w:t xml:space="preserve">a) </w:t></w:r><w:r><w:t>The solution for he sub did meet the goal to find the quantity of each color of socks</w:t></w:r></w:p>'
$inner50 += '<w:p/>'
$inner50 += '<w:p><w:r><w:t>b)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>his</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> solution to find the number in each sock will work for all cases and all colors.</w:t></w:r></w:p>'
$inner50 += '<w:p/>'
$inner50 += '<w:p><w:r><w:t xml:space="preserve">5. </w:t></w:r></w:p>'
$inner50 += '<w:p/>'
$inner50 += '<w:p><w:r><w:t xml:space="preserve">a) You would pick five times to get the first matching pair of socks witch is half the number of the color of socks with the most individual socks.  You would have to pick ten times to get a total of each matching color witch is half the number of sock within each color </w:t></w:r></w:p>'
$p50.InsertXML((New-XmlFrag $inner50))
